$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the style already used by the
# other header cells (e.g. G1: bold font, border, centered/top alignment).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill column H (rows 2-42) with the "Save" flag (1 when the pitching
# outing met the save criteria, 0 otherwise), mirroring the rest of the
# data which has no special cell style.
$saveVals = @(0,0,0,1,1,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,1,1,1,0,0,0,0,1,0,1)
for ($i = 0; $i -lt $saveVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveVals[$i]
}
